$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15:A26").NumberFormat = "@"

# Row 15
$ws.Range("A15").Value = "2022"
$ws.Range("B15").Value = "JANUARI"
$ws.Range("C15").Value = 44571
$ws.Range("D15").Value = "MAAS LIMBURG"
$ws.Range("E15").Value = "Denderbekken"
$ws.Range("F15").Value = "Niet geklasseerd"
$ws.Range("H15").Value = "Oost-Vlaanderen"
$ws.Range("I15").Value = "Aalst"
$ws.Range("J15").Value = "Waterloop"
$ws.Range("K15").Value = "Vuilzakken"
$ws.Range("L15").Value = 137299
$ws.Range("M15").Value = 4.0444033023359198
$ws.Range("N15").Value = 50.944339968671798
$ws.Range("O15").Value = "Regulier"
$ws.Range("P15").Value = 36526
$ws.Range("Q15").Value = 2958465
$ws.Range("R15").Value = 535509
$ws.Range("S15").Value = "Niet van toepassing"

# Row 16
$ws.Range("A16").Value = "2022"
$ws.Range("B16").Value = "JANUARI"
$ws.Range("C16").Value = 44571
$ws.Range("D16").Value = "MAAS LIMBURG"
$ws.Range("E16").Value = "Denderbekken"
$ws.Range("F16").Value = "Niet geklasseerd"
$ws.Range("H16").Value = "Oost-Vlaanderen"
$ws.Range("I16").Value = "Aalst"
$ws.Range("J16").Value = "Waterloop"
$ws.Range("K16").Value = "Vuilzakken"
$ws.Range("L16").Value = 137309
$ws.Range("M16").Value = 4.0444710281176102
$ws.Range("N16").Value = 50.944348015299298
$ws.Range("O16").Value = "Regulier"
$ws.Range("P16").Value = 36526
$ws.Range("Q16").Value = 2958465
$ws.Range("R16").Value = 535525
$ws.Range("S16").Value = "Niet van toepassing"

# Row 17
$ws.Range("A17").Value = "2022"
$ws.Range("B17").Value = "JANUARI"
$ws.Range("C17").Value = 44571
$ws.Range("D17").Value = "LEIE & BOVENSCHELDE"
$ws.Range("E17").Value = "Bovenscheldebekken"
$ws.Range("F17").Value = "Onbevaarbaar cat. 2"
$ws.Range("H17").Value = "West-Vlaanderen"
$ws.Range("I17").Value = "Spiere-Helkijn"
$ws.Range("J17").Value = "Waterloop"
$ws.Range("K17").Value = "Muskusrat"
$ws.Range("L17").Value = 137338
$ws.Range("M17").Value = 3.3398925641114898
$ws.Range("N17").Value = 50.73213179303
$ws.Range("O17").Value = "Regulier"
$ws.Range("P17").Value = 36526
$ws.Range("Q17").Value = 2958465
$ws.Range("R17").Value = 535572
$ws.Range("S17").Value = "Niet van toepassing"

# Row 18
$ws.Range("A18").Value = "2022"
$ws.Range("B18").Value = "JANUARI"
$ws.Range("C18").Value = 44571
$ws.Range("D18").Value = "MAAS LIMBURG"
$ws.Range("E18").Value = "Denderbekken"
$ws.Range("F18").Value = "Niet geklasseerd"
$ws.Range("H18").Value = "Oost-Vlaanderen"
$ws.Range("I18").Value = "Aalst"
$ws.Range("J18").Value = "Waterloop"
$ws.Range("K18").Value = "Vuilzakken"
$ws.Range("L18").Value = 137398
$ws.Range("M18").Value = 4.0444626462139297
$ws.Range("N18").Value = 50.944299064981898
$ws.Range("O18").Value = "Regulier"
$ws.Range("P18").Value = 36526
$ws.Range("Q18").Value = 2958465
$ws.Range("R18").Value = 535662
$ws.Range("S18").Value = "Niet van toepassing"

# Row 19
$ws.Range("A19").Value = "2022"
$ws.Range("B19").Value = "JANUARI"
$ws.Range("C19").Value = 44572
$ws.Range("D19").Value = "BRUGSE POL-GENTSE KAN"
$ws.Range("E19").Value = "Bekken van de Gentse Kanalen"
$ws.Range("F19").Value = "Bevaarbaar"
$ws.Range("G19").Value = "LEOPOLDKANAAL"
$ws.Range("H19").Value = "Oost-Vlaanderen"
$ws.Range("I19").Value = "Sint-Laureins"
$ws.Range("J19").Value = "Waterloop"
$ws.Range("K19").Value = "Muskusrat"
$ws.Range("L19").Value = 137471
$ws.Range("M19").Value = 3.58982886164361
$ws.Range("N19").Value = 51.260769908027001
$ws.Range("O19").Value = "Regulier"
$ws.Range("P19").Value = 36526
$ws.Range("Q19").Value = 2958465
$ws.Range("R19").Value = 535878
$ws.Range("S19").Value = "Niet van toepassing"

# Row 20
$ws.Range("A20").Value = "2022"
$ws.Range("B20").Value = "JANUARI"
$ws.Range("C20").Value = 44573
$ws.Range("D20").Value = "DIJLE EN ZENNE"
$ws.Range("E20").Value = "Dijle- en Zennebekken"
$ws.Range("F20").Value = "Bevaarbaar"
$ws.Range("G20").Value = "KANAAL LEUVEN-DIJLE - KANAAL VAN LEUVEN NAAR DE DIJLE"
$ws.Range("H20").Value = "Vlaams-Brabant"
$ws.Range("I20").Value = "Kampenhout"
$ws.Range("J20").Value = "Waterloop"
$ws.Range("K20").Value = "Grote waternavel"
$ws.Range("L20").Value = 137607
$ws.Range("M20").Value = 4.6049863609663397
$ws.Range("N20").Value = 50.954578464006097
$ws.Range("O20").Value = "ANB-Waterplanten"
$ws.Range("P20").Value = 43801
$ws.Range("Q20").Value = 44727
$ws.Range("R20").Value = 536201
$ws.Range("S20").Value = "Waargenomen        "

# Row 21
$ws.Range("A21").Value = "2022"
$ws.Range("B21").Value = "JANUARI"
$ws.Range("C21").Value = 44573
$ws.Range("D21").Value = "BRUGSE POL-GENTSE KAN"
$ws.Range("E21").Value = "Bekken van de Brugse Polders"
$ws.Range("F21").Value = "Polder of wateringgracht"
$ws.Range("H21").Value = "West-Vlaanderen"
$ws.Range("I21").Value = "Knokke-Heist"
$ws.Range("J21").Value = "Waterloop"
$ws.Range("K21").Value = "Bruine rat"
$ws.Range("L21").Value = 137691
$ws.Range("M21").Value = 3.3678554329488
$ws.Range("N21").Value = 51.334779099957899
$ws.Range("O21").Value = "Regulier"
$ws.Range("P21").Value = 36526
$ws.Range("Q21").Value = 2958465
$ws.Range("R21").Value = 536372
$ws.Range("S21").Value = "Niet van toepassing"

# Row 22
$ws.Range("A22").Value = "2022"
$ws.Range("B22").Value = "JANUARI"
$ws.Range("C22").Value = 44573
$ws.Range("D22").Value = "DEMER"
$ws.Range("E22").Value = "Demerbekken"
$ws.Range("F22").Value = "Onbevaarbaar cat. 2"
$ws.Range("G22").Value = "ZUTENDAALBEEK"
$ws.Range("H22").Value = "Limburg"
$ws.Range("I22").Value = "Zutendaal"
$ws.Range("J22").Value = "Waterloop"
$ws.Range("K22").Value = "Bever"
$ws.Range("L22").Value = 137759
$ws.Range("M22").Value = 5.5332117845542301
$ws.Range("N22").Value = 50.908417141194199
$ws.Range("O22").Value = "Regulier"
$ws.Range("P22").Value = 36526
$ws.Range("Q22").Value = 2958465
$ws.Range("R22").Value = 536538
$ws.Range("S22").Value = "Niet van toepassing"

# Row 23
$ws.Range("A23").Value = "2022"
$ws.Range("B23").Value = "JANUARI"
$ws.Range("C23").Value = 44575
$ws.Range("D23").Value = "BRUGSE POL-GENTSE KAN"
$ws.Range("E23").Value = "Bekken van de Brugse Polders"
$ws.Range("F23").Value = "Onbevaarbaar cat. 1"
$ws.Range("G23").Value = "ISABELLAVAART - KLEINE GEULE"
$ws.Range("H23").Value = "West-Vlaanderen"
$ws.Range("I23").Value = "Knokke-Heist"
$ws.Range("J23").Value = "Waterloop"
$ws.Range("K23").Value = "Woelrat"
$ws.Range("L23").Value = 137963
$ws.Range("M23").Value = 3.3684540685089801
$ws.Range("N23").Value = 51.335805883157597
$ws.Range("O23").Value = "Regulier"
$ws.Range("P23").Value = 36526
$ws.Range("Q23").Value = 2958465
$ws.Range("R23").Value = 536884
$ws.Range("S23").Value = "Niet van toepassing"

# Row 24
$ws.Range("A24").Value = "2022"
$ws.Range("B24").Value = "JANUARI"
$ws.Range("C24").Value = 44575
$ws.Range("D24").Value = "LEIE & BOVENSCHELDE"
$ws.Range("E24").Value = "Bovenscheldebekken"
$ws.Range("F24").Value = "Onbevaarbaar cat. 2"
$ws.Range("G24").Value = "RIJTGRACHT - PACHTBEEK - OLIEBERGBEEK"
$ws.Range("H24").Value = "West-Vlaanderen"
$ws.Range("I24").Value = "Avelgem"
$ws.Range("J24").Value = "Natuurgebied"
$ws.Range("K24").Value = "Muskusrat"
$ws.Range("L24").Value = 138035
$ws.Range("M24").Value = 3.4821835931946898
$ws.Range("N24").Value = 50.783272301898698
$ws.Range("O24").Value = "Regulier"
$ws.Range("P24").Value = 36526
$ws.Range("Q24").Value = 2958465
$ws.Range("R24").Value = 536986
$ws.Range("S24").Value = "Niet van toepassing"

# Row 25
$ws.Range("A25").Value = "2022"
$ws.Range("B25").Value = "JANUARI"
$ws.Range("C25").Value = 44579
$ws.Range("D25").Value = "LEIE & BOVENSCHELDE"
$ws.Range("E25").Value = "Leiebekken"
$ws.Range("F25").Value = "Onbevaarbaar cat. 2"
$ws.Range("G25").Value = "TOLBEEK"
$ws.Range("H25").Value = "West-Vlaanderen"
$ws.Range("I25").Value = "Wevelgem"
$ws.Range("J25").Value = "Waterloop"
$ws.Range("K25").Value = "Grote waternavel"
$ws.Range("L25").Value = 138443
$ws.Range("M25").Value = 3.2085166176269402
$ws.Range("N25").Value = 50.809415627087098
$ws.Range("O25").Value = "Regulier"
$ws.Range("P25").Value = 36526
$ws.Range("Q25").Value = 2958465
$ws.Range("R25").Value = 537859
$ws.Range("S25").Value = "Waargenomen        "

# Row 26
$ws.Range("A26").Value = "2022"
$ws.Range("B26").Value = "JANUARI"
$ws.Range("C26").Value = 44579
$ws.Range("D26").Value = "LEIE & BOVENSCHELDE"
$ws.Range("E26").Value = "Leiebekken"
$ws.Range("F26").Value = "Onbevaarbaar cat. 2"
$ws.Range("G26").Value = "TOLBEEK"
$ws.Range("H26").Value = "West-Vlaanderen"
$ws.Range("I26").Value = "Wevelgem"
$ws.Range("J26").Value = "Waterloop"
$ws.Range("K26").Value = "Grote waternavel"
$ws.Range("L26").Value = 138446
$ws.Range("M26").Value = 3.2083872010342702
$ws.Range("N26").Value = 50.809368185512298
$ws.Range("O26").Value = "Regulier"
$ws.Range("P26").Value = 36526
$ws.Range("Q26").Value = 2958465
$ws.Range("R26").Value = 537862
$ws.Range("S26").Value = "Waargenomen        "

# Update view: scroll to top-left A3, select A5:S26
$ws.Range("A5:S26").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
